$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo: reorder phrase in do_not_use for row 6 (The Ordinary "Buffet" + Copper Peptides 1%) ---
$ws.Range("H6").Value = "with Direct Acids Direct Vitamin C EUK Retinoids Strong Antioxidants Resveratrol and Ferulic Acid Salicylic Acid"

# --- Reorder product rows 9-13 ---
# Target order: Granactive Retinoid 5% in Squalane, Marine Hyaluronics, Azelaic Acid Suspension 10%,
#               Granactive Retinoid 2% Emulsion, Amino Acids + B5
# Row 9: The Ordinary Granactive Retinoid 5% in Squalane
$ws.Range("A9").Value = "The Ordinary Granactive Retinoid 5% in Squalane"
$ws.Range("B9").Value = "Textural Irregularities, Signs of Aging, Uneven Skin Tone"
$ws.Range("C9").Value = "Dry Skin"
$ws.Range("D9").Value = "Anhydrous Serum"
$ws.Range("E9").Value = "Granactive Retinoid 5% in Squalane is a light, oil-like serum that targets visible signs of aging through a highly-advanced retinoid active, one of the best all-around ingredients for skin. This next-generation retinoid improves the appearance of fine and dynamic lines that arise from a loss of collagen and elastin while evening skin texture and tone. It’s further supported with squalane, a moisturizing agent naturally found in the skin. That means it not only delivers comparable visible results to retinol, it’s less likely to cause irritation.Granactive Retinoid™ is a trademark of Grant Industries. Neither DECIEM nor The Ordinary is affiliated with Grant Industries.Note: Once opened, this formula should be refrigerated and used within a three-month period.Caution: Retinoids can make the skin more sensitive to UV radiation. Sun protection is particularly important when using retinoids. This product must not be used in conjunction with other retinoids including retinol or retinoic acid. This product is not a treatment for acne.Note: When pregnant or breastfeeding, it is recommended to avoid any skincare products containing retinoids such as formulations with Granactive Retinoid or Retinol."
$ws.Range("F9").Value = "Use in PM"
$ws.Range("G9").Value = "3 months after opening."
$ws.Range("H9").Value = "with Copper Peptides Direct Acids Direct Vitamin C Retinoids"

# Row 10: The Ordinary Marine Hyaluronics
$ws.Range("A10").Value = "The Ordinary Marine Hyaluronics"
$ws.Range("B10").Value = "Dryness"
$ws.Range("C10").Value = "All Skin Types"
$ws.Range("D10").Value = "Water-based Serum"
$ws.Range("E10").Value = "Marine Hyaluronics offers an exceptionally lightweight serum, acting as a hyaluronic acid-alternative, which directs water where you need it most. By combining exopolysaccharides from skin-friendly marine bacteria, Hawaiian red algae, glycoproteins from Antarctic marine sources, micro-filtered blue-green algae, and several amino acids, this water-based formula targets hydration below the skin surface, resulting in a softer, plumper complexion. "
$ws.Range("F10").Value = "Use in AM Use in PM"
$ws.Range("G10").Value = "12 months after opening."
$ws.Range("H10").Value = ""

# Row 11: The Ordinary Azelaic Acid Suspension 10%
$ws.Range("A11").Value = "The Ordinary Azelaic Acid Suspension 10%"
$ws.Range("B11").Value = "Textural Irregularities, Dullness, Uneven Skin Tone, Look of Redness"
$ws.Range("C11").Value = "All Skin Types"
$ws.Range("D11").Value = "Suspension"
$ws.Range("E11").Value = "Azelaic Acid Suspension 10% is a cream-like formula that brightens skin tone and visibly improves skin texture due to a high concentration of azelaic acid—a natural and effective antioxidant found in grains. It has a lower irritation potential than other direct acids, which means it’s gentle enough for daily use as part of your skincare regimen."
$ws.Range("F11").Value = "Use in AM Use in PM"
$ws.Range("G11").Value = "12 months after opening."
$ws.Range("H11").Value = "with Copper Peptides Direct Acids Direct Vitamin C EUK Niacinamide Powder Peptides Retinoids"

# Row 12: The Ordinary Granactive Retinoid 2% Emulsion
$ws.Range("A12").Value = "The Ordinary Granactive Retinoid 2% Emulsion"
$ws.Range("B12").Value = "Textural Irregularities, Signs of Aging, Uneven Skin Tone"
$ws.Range("C12").Value = "All Skin Types"
$ws.Range("D12").Value = "Emulsion"
$ws.Range("E12").Value = "Granactive Retinoid 2% Emulsion is a creamy serum that targets visible signs of aging through a highly-advanced retinoid active, one of the best all-around ingredients for skin. This next-generation retinoid improves the appearance of fine and dynamic lines that arise from a loss of collagen and elastin, while evening skin texture and tone. What’s more, this advanced ingredient is not only proven to deliver comparable visible results to retinol; it’s less likely to cause irritation.Granactive Retinoid™ is a trademark of Grant Industries. Neither DECIEM nor The Ordinary is affiliated with Grant Industries.Note: Once opened, this formula should be refrigerated and used within a three-month period.Caution: Retinoids can make the skin more sensitive to UV radiation. Sun protection is particularly important when using retinoids. This product must not be used in conjunction with other retinoids including retinol or retinoic acid. This product is not a treatment for acne.Note: When pregnant or breastfeeding, it is recommended to avoid any skincare products containing retinoids such as formulations with Granactive Retinoid or Retinol."
$ws.Range("F12").Value = "Use in PM"
$ws.Range("G12").Value = "3 months after opening."
$ws.Range("H12").Value = "with Copper Peptides Direct Acids Direct Vitamin C Retinoids"

# Row 13: The Ordinary Amino Acids + B5
$ws.Range("A13").Value = "The Ordinary Amino Acids + B5"
$ws.Range("B13").Value = "Dryness"
$ws.Range("C13").Value = "All Skin Types"
$ws.Range("D13").Value = "Water-based Serum"
$ws.Range("E13").Value = "Amino Acids + B5 is an ultra-thin serum that works with your skin to support its natural hydration barrier. By including a concentrated 17% (by weight) solution of amino acids and amino acid derivatives that mimic your skin’s natural moisturizing factors, this water-based formula keeps the outer layer of your skin protected and well-hydrated without feeling greasy. Plus, it uses 5% (by weight) pro-vitamin B5 to provide surface and below-surface hydration, giving way to softer, smoother skin."
$ws.Range("F13").Value = "Use in AM Use in PM"
$ws.Range("G13").Value = "12 months after opening."
$ws.Range("H13").Value = ""

# --- Reorder product rows 25-27 ---
# Target order: Alpha Arbutin 2% + HA, 100% Plant-Derived Squalane, Aloe 2% + NAG 2% Solution
# Row 25: The Ordinary Alpha Arbutin 2% + HA
$ws.Range("A25").Value = "The Ordinary Alpha Arbutin 2% + HA"
$ws.Range("B25").Value = "Uneven Skin Tone, Dryness"
$ws.Range("C25").Value = "All Skin Types"
$ws.Range("D25").Value = "Water-based Serum"
$ws.Range("E25").Value = "Alpha Arbutin 2% + HA is a water-based serum specifically designed to target uneven skin tone and visibly improve pigmentation. It combines a high concentration of purified alpha arbutin, a well-known skin-brightening ingredient, with hyaluronic acid. Purified alpha-arbutin works to target dark spots and uneven skin tone. Meanwhile, hyaluronic acid helps to support product absorption into the skin.Note: Alpha Arbutin is extremely sensitive to degradation in the presence of water if the pH of the formulation is not ideal. The pH of this formula has been shown to be the most suitable pH to minimize degradation of Alpha Arbutin.We are aware of the potential colour changes of our Alpha Arbutin product. This is a natural occurrence, where certain ingredients in the product can change colour when exposed to sunlight or increased temperature. Testing has shown that the efficacy and safety profile of the active ingredient remains stable."
$ws.Range("F25").Value = "Use in AM Use in PM"
$ws.Range("G25").Value = "12 months after opening."
$ws.Range("H25").Value = ""

# Row 26: The Ordinary 100% Plant-Derived Squalane
$ws.Range("A26").Value = "The Ordinary 100% Plant-Derived Squalane"
$ws.Range("B26").Value = "Dryness, Hair"
$ws.Range("C26").Value = "All Skin Types"
$ws.Range("D26").Value = "Anhydrous Serum"
$ws.Range("E26").Value = "100% Plant-Derived Squalane hydrates your skin while supporting its natural moisture barrier. Squalane is an exceptional hydrator found naturally in the skin, and this formula uses 100% plant-derived squalane derived from sugar cane for a non-comedogenic solution that enhances surface-level hydration.Our 100% Plant-Derived Squalane formula can also be used in hair to increase heat protection, add shine, and reduce breakage."
$ws.Range("F26").Value = "Use in AM Use in PM"
$ws.Range("G26").Value = "6 months after opening."
$ws.Range("H26").Value = ""

# Row 27: Aloe 2% + NAG 2% Solution
$ws.Range("A27").Value = "Aloe 2% + NAG 2% Solution"
$ws.Range("B27").Value = "Textural Irregularities, Uneven Skin Tone, Look of Redness"
$ws.Range("C27").Value = "Dry Skin"
$ws.Range("D27").Value = "Water-based Serum"
$ws.Range("E27").Value = "Aloe 2% + NAG 2% Solution is a lightweight soothing serum formulated for blemish-prone skin. It reduces the appearance of post-acne marks through the inclusion of ingredients that target the look and feel of both uneven skin tone and texture, while also helping to reduce the appearance of pores, and strengthening skin barrier.This formulation incorporates aloe barbadensis leaf juice powder for hydration. It also contains N-acetyl glucosamine (NAG), a biotechnology-derived ingredient that has been shown to effectively target the appearance of uneven skin tone. It is further supported by peptide technology palmitoyl pentapeptide-4 which targets the feel of textural irregularities, while offering effective barrier support.Testing ShowsReduces the appearance of post-acne marksReduces the look of skin redness caused by irritationReduces the appearance of poresImproves the appearance of uneven skin tone and the feel of uneven skin textureStrengthens and maintains skin barrier functionBoosts skin hydration*Consumer testing on 31 subjects after using product 2x/day for 8 weeks."
$ws.Range("F27").Value = "Use in AM Use in PM"
$ws.Range("G27").Value = "12 months after opening."
$ws.Range("H27").Value = "with Direct Acids Direct Vitamin C Resveratrol and Ferulic Acid Salicylic Acid"

